# Iteration v0.7.2 -- Additional URL cleanup
#
# Updates the "Results" sheet of the workbook with the metrics for the new
# 0.7.2 iteration (stored in column J) and makes "Results" the active /
# selected sheet+cell (it was "Steps" before).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Results")

# Label the new iteration column (merged J8:K8 header) with its version tag.
$ws.Range("J8").Value = "0.7.2"

# New iteration metrics for 0.7.2.
$ws.Range("J10").Value = 0.922979797979798    # Accuracy
$ws.Range("J11").Value = 0.0547877591312932   # FPR
$ws.Range("J12").Value = 0.862612612612613    # F1

# Fill in the previously-empty Test columns (I12, K12) for the 0.6 iteration.
$ws.Range("I12").Value = 0.870070989112258
$ws.Range("K12").Value = 0.871422594308083

# Make "Results" the active sheet with K12 selected (previously "Steps" was
# the active sheet with B21 selected there, and H12 was selected on Results).
[void]$ws.Activate()
[void]$ws.Range("K12").Select()
